$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{Row=1; Col=1; New="49÷5="},
    @{Row=1; Col=2; New="64÷7="},
    @{Row=1; Col=3; New="63÷5="},
    @{Row=1; Col=4; New="61÷8="},
    @{Row=1; Col=5; New="84÷7="},
    @{Row=5; Col=1; New="97÷4="},
    @{Row=5; Col=2; New="79÷3="},
    @{Row=5; Col=3; New="77÷8="},
    @{Row=5; Col=4; New="21÷9="},
    @{Row=5; Col=5; New="57÷7="},
    @{Row=9; Col=1; New="84÷2="},
    @{Row=9; Col=2; New="90÷8="},
    @{Row=9; Col=3; New="48÷6="},
    @{Row=9; Col=4; New="32÷2="},
    @{Row=9; Col=5; New="98÷5="},
    @{Row=13; Col=1; New="19÷4="},
    @{Row=13; Col=2; New="41÷3="},
    @{Row=13; Col=3; New="65÷5="},
    @{Row=13; Col=4; New="30÷4="},
    @{Row=13; Col=5; New="43÷2="},
    @{Row=17; Col=1; New="34÷2="},
    @{Row=17; Col=2; New="28÷4="},
    @{Row=17; Col=3; New="53÷9="},
    @{Row=17; Col=4; New="54÷6="},
    @{Row=17; Col=5; New="22÷7="}
)

foreach ($c in $changes) {
    $cell = $t.Cell($c.Row, $c.Col)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $c.New
}
